$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17287266254425
$ws.Range("B1").Value = 2.388129949569702
$ws.Range("D1").Value = 2.367826700210571
$ws.Range("E1").Value = 1.20968770980835
